$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "mu" column entirely; this shifts isSelected/bandwidth/transRate/uploadTime/totalTime
# left by one column (J->I, K->J, L->K, M->L, N->M) and updates the used range to A1:M21.
$ws.Columns("I").Delete()

# Refresh the simulation results for the new "Eligible for Experiment" round (C2:M21).
$rows = 2..21
$newValues = @{
    2 = @(70.95848971701832, 79561592, 1126927007.043005, 0.004220205241383404, 25.04851854875356, 20, $true, 71428.57142857143, 101898.3714322749, 19.47226410108782, 44.52078264984138)
    3 = @(175.9556230791904, 79873912, 232371787.401807, 0.003491682868916071, 302.4089996224546, 0, $false, 71428.57142857143, 90073.79794100023, 22.02851489952358, 324.4375145219782)
    4 = @(118.5748232692446, 79817080, 1099128883.62274, 0.003200154175598993, 43.0536231732567, 20, $true, 71428.57142857143, 84933.77474169277, 23.36163682862889, 66.4152600018856)
    5 = @(80.7821316518251, 79235960, 833180992.2729816, 0.0029241020420135, 38.41212060549256, 20, $true, 71428.57142857143, 79818.30177382777, 24.85886013488966, 63.27098074038223)
    6 = @(80.74375796642178, 79761784, 895042835.8966057, 0.00430839303353433, 35.97741875568727, 20, $true, 71428.57142857143, 103241.8734391192, 19.21886860344569, 55.19628735913295)
    7 = @(111.7396685830363, 80025976, 1405892873.14731, 0.002292505142216274, 31.80212449706709, 20, $true, 71428.57142857143, 67060.84892685739, 29.58793441705664, 61.39005891412373)
    8 = @(123.9098299632842, 79814520, 1235208275.756846, 0.0009087328629546882, 40.03293936701237, 9, $false, 71428.57142857143, 31950.53194814984, 62.10200203301775, 102.1349414000301)
    9 = @(110.2705246966666, 79492984, 967188858.4782858, 0.002054846491575216, 45.31551919019817, 15, $true, 71428.57142857143, 61819.59556699799, 32.09648949983181, 77.41200869002998)
    10 = @(109.2942173926688, 80046456, 1828838693.402923, 0.002065237360215993, 23.91849755567599, 20, $true, 71428.57142857143, 62054.41266155088, 31.97503472995422, 55.89353228563021)
    11 = @(153.218164522722, 79497592, 788234814.5322708, 0.002172038545255918, 77.26425492538016, 1, $false, 71428.57142857143, 64437.43442166985, 30.79253570239491, 108.0567906277751)
    12 = @(52.69345785456528, 79762808, 1304380947.56723, 0.003557886963972279, 16.11100717757591, 20, $true, 71428.57142857143, 91206.18312816262, 21.75501629326841, 37.86602347084433)
    13 = @(185.9418346898805, 80041848, 983761988.3588175, 0.002715096879788233, 75.64394764793487, 2, $false, 71428.57142857143, 75768.99705409686, 26.18738636045749, 101.8313340083924)
    14 = @(197.075549015551, 79997816, 1645456896.429853, 0.001254936285397894, 47.90649193683429, 9, $false, 71428.57142857143, 41918.59522046038, 47.33441064913168, 95.24090258596597)
    15 = @(72.50100601321395, 79731064, 1371947743.31511, 0.003943618243461979, 21.06706461186349, 20, $true, 71428.57142857143, 97567.35201726464, 20.33663883436024, 41.40370344622373)
    16 = @(96.3739219971601, 79832952, 789119212.298843, 0.002360910960038808, 48.74938139218275, 15, $true, 71428.57142857143, 68521.32948914779, 28.95728986569431, 77.70667125787706)
    17 = @(132.3194100858378, 79816568, 1734880216.672365, 0.003627136141331431, 30.4380702809947, 20, $true, 71428.57142857143, 92377.48661516073, 21.47917282341778, 51.91724310441248)
    18 = @(176.1311528522601, 79846264, 1118402169.629412, 0.004926720339789661, 62.87279706335826, 14, $true, 71428.57142857143, 112199.7970280936, 17.68445266886874, 80.55724973222699)
    19 = @(102.2661244817619, 79854968, 1541271038.75871, 0.004131978223816866, 26.49260867365715, 20, $true, 71428.57142857143, 100536.5124890668, 19.73603371427647, 46.22864238793362)
    20 = @(165.4740402215803, 80000888, 1386432394.135247, 0.002262925582568383, 47.74149181262842, 15, $true, 71428.57142857143, 66422.85091959659, 29.87212943331537, 77.61362124594379)
    21 = @(71.22950717745928, 80230776, 427332063.0859565, 0.003963270018734173, 66.86601742069158, 9, $false, 71428.57142857143, 97881.1550043938, 20.27144040046249, 87.13745782115407)
}

foreach ($r in $rows) {
    $rowData = $newValues[$r]
    $arr = New-Object "object[,]" 1,11
    for ($c = 0; $c -lt 11; $c++) {
        $arr[0, $c] = $rowData[$c]
    }
    $rangeStr = "C" + $r + ":M" + $r
    $ws.Range($rangeStr).Value = $arr
}